$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.967.61'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.555.40'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = "'207.01"
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('D6').Value = "'0.486"
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').Value = "'22.15"
$ws.Range('E8').Value = '  +3.93%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').Value = "'0.0587"
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('D11').Value = "'0.0858"
$ws.Range('D12').Value = '1.777.35'
$ws.Range('D13').Value = '1.555.58'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').Value = "'3.75"
$ws.Range('E14').Value = '  +1.27%  '
$ws.Range('D15').Value = "'0.519"
$ws.Range('E15').Value = '  +1.51%  '
$ws.Range('D16').Value = '26.965.15'
$ws.Range('E16').Value = '  +0.34%  '
$ws.Range('D17').Value = "'61.68"
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('D18').Value = "'218.01"
$ws.Range('E18').Value = '  +2.03%  '
$ws.Range('E19').Value = '  +2.24%  '
$ws.Range('D20').Value = "'7.30"
$ws.Range('E21').Value = '  -0.15%  '
$ws.Range('E22').Value = '  +1.24%  '
$ws.Range('E23').Value = '  +0.48%  '
$ws.Range('D24').Value = "'1.94"
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('D25').Value = "'154.54"
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  +0.79%  '
$ws.Range('E28').Value = '  +0.97%  '
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('E30').Value = '  +2.11%  '
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('E32').Value = '  +0.47%  '
$ws.Range('D33').Value = '1.426.73'
$ws.Range('E33').Value = '  +4.78%  '
$ws.Range('E34').Value = '  +4.48%  '
$ws.Range('E35').Value = '  +3.60%  '
$ws.Range('E36').Value = '  +1.82%  '
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('E38').Value = '  +0.74%  '
$ws.Range('D39').Value = "'0.522"
$ws.Range('E39').Value = '  +0.77%  '
$ws.Range('E40').Value = '  +0.66%  '
$ws.Range('D41').Value = "'5.79"
$ws.Range('E41').Value = '  +3.67%  '
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('E43').Value = '  +4.47%  '
$ws.Range('E44').Value = '  +0.37%  '
$ws.Range('D45').Value = "'64.37"
$ws.Range('E45').Value = '  +1.31%  '
$ws.Range('E46').Value = '  +1.28%  '
$ws.Range('D47').Value = '1.691.16'
$ws.Range('E47').Value = '  +0.54%  '
$ws.Range('D48').Value = "'87.83"
$ws.Range('E48').Value = '  +2.02%  '
$ws.Range('D49').Value = "'0.0521"
$ws.Range('E49').Value = '  +2.53%  '
$ws.Range('D50').Value = '0.0₆0100'
$ws.Range('E50').Value = '  +3.78%  '
$ws.Range('E51').Value = '  +0.79%  '
